$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Iterative calculation delta: 0.001 -> 0.0001 ---
$excel.Iteration = $false
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# --- Remove the two trailing empty styled rows (old rows 4 & 5) ---
$ws.Rows("4:5").Delete()

# --- Add the new rows of data (Path / Text pairs) ---
$ws.Range("A2").Value = "/Users/lanhdang/Projects/mmo/add-text-to-videos/video1.mp4"
$ws.Range("B2").Value = "Thuốc thông minh, uống một liều duy nhất"

$ws.Range("A3").Value = "/Users/lanhdang/Projects/mmo/add-text-to-videos/video3.mp4"
$ws.Range("B3").Value = "Thuốc thử độ ngu, uống càng nhiều càng cho kết quả chính xác"

$ws.Range("A4").Value = "/Users/lanhdang/Projects/mmo/add-text-to-videos/video12.mp4"
$ws.Range("B4").Value = "Thuốc làm giàu, uống vào giàu ngay!!!!"

# --- Column widths ---
$ws.Columns("A:A").ColumnWidth = 49.8
$ws.Columns("B:B").ColumnWidth = 52.3

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 169
[void]$ws.Range("A5:B8").Select()
